$d = $word.ActiveDocument

# 1. Insert a new "Compact" list item ("Full Stack Web Development with Flask
#    (Dec 2020)") immediately before the existing "Python Decorators (Dec 2020)"
#    entry under the "LinkedIn Learning" heading. Locate the target paragraph,
#    then add a new paragraph before it (inherits the same paragraph style /
#    numbering properties) and fill in its text.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Python Decorators (Dec 2020)*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($targetIndex)
    $newPara.Range.InsertBefore("Full Stack Web Development with Flask (Dec 2020)")
}

# 2. Bump the "last update" footer date.
$d.Content.Find.Execute("Last update: December 05, 2020", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Last update: December 09, 2020", 2)
